$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

# Sample data row
$ws.Range("A2").Value = 41000000005
$ws.Range("B2").Value = "Headphone (HP123)"
$ws.Range("C2").Value = "HN @ DM Alexandria"
$ws.Range("D2").Value = 1000
$ws.Range("D2").NumberFormat = "`"$`"#,##0;[Red]\-`"$`"#,##0"

# Size the header/data columns to fit their contents
$ws.Range("A1:C2").EntireColumn.AutoFit()

# Leave the selection on the data row, matching the saved view
$ws.Range("A2").Select() | Out-Null
